$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.9649872183799744
$ws.Range("B1").Value = 4.451056480407715
$ws.Range("C1").Value = 1.755985140800476
$ws.Range("D1").Value = 0.9119553565979004
$ws.Range("E1").Value = 0.9552488923072815
